$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Phase 1: propagate existing formats to their new homes BEFORE the donor
# cells' own formats/values are changed.
# ---------------------------------------------------------------------------

# 1a. D1 currently carries numFmt #,##0 / border / (theme) font / right-align.
#     That is the format the new "empty numeric" cells in rows 3-5 need.
$ws.Range("D1").Copy()
$ws.Range("E3:G3").PasteSpecial(-4122)
$ws.Range("E4:G4").PasteSpecial(-4122)
$ws.Range("E5:G5").PasteSpecial(-4122)

# 1b. A1 currently carries numFmt #,##0 / right-align with NO border/font -
#     exactly what D2 needs once its value is cleared.
$ws.Range("A1").Copy()
$ws.Range("D2").PasteSpecial(-4122)

# 1c. B1 currently carries the plain General / general-align format that
#     column A (all rows) and the new B5:C5 cells need.
$ws.Range("B1").Copy()
$ws.Range("A1:A5").PasteSpecial(-4122)
$ws.Range("B5:C5").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Phase 2: build the new "black font" numeric style on D1, then fan it out.
# ---------------------------------------------------------------------------
$ws.Range("D1").Font.Color = 0
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Phase 3: set the final cell values.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = $null
$ws.Range("B1").Value = "create"
$ws.Range("C1").Value = "python.Array"
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 3

$ws.Range("A2").Value = $null
$ws.Range("B2").Value = "__len__"
$ws.Range("C2").Value = "A1"
$ws.Range("D2").Value = $null
$ws.Range("E2").Value = $null
$ws.Range("F2").Value = $null
$ws.Range("G2").Value = $null

$ws.Range("A3").Value = "wow"
$ws.Range("B3").Value = "create"
$ws.Range("C3").Value = "Calculator"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = $null
$ws.Range("F3").Value = $null
$ws.Range("G3").Value = $null

$ws.Range("A4").Value = $null
$ws.Range("B4").Value = "addme"
$ws.Range("C4").Value = "A2"
$ws.Range("D4").Value = 18
$ws.Range("E4").Value = $null
$ws.Range("F4").Value = $null
$ws.Range("G4").Value = $null

$ws.Range("A5").Value = $null
$ws.Range("B5").Value = "subme"
$ws.Range("C5").Value = "A2"
$ws.Range("D5").Value = 16
$ws.Range("E5").Value = $null
$ws.Range("F5").Value = $null
$ws.Range("G5").Value = $null

# ---------------------------------------------------------------------------
# Phase 4: row heights.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(3).RowHeight = 18.75
$ws.Rows.Item(4).RowHeight = 18.75
$ws.Rows.Item(5).RowHeight = 18.75
